# Remove the unneeded "import java.util.Scanner;" line (and the blank
# line that followed it) from the first two Java code examples
# (slides 26 and 27), since those examples no longer use Scanner.

$p = $ppt.ActivePresentation

$slideIndexes = @(26, 27)

foreach ($idx in $slideIndexes) {
    $s = $p.Slides.Item($idx)
    $shp = $s.Shapes.Item(1)
    $tr = $shp.TextFrame.TextRange

    # Paragraph 1 is "import java.util.Scanner;" and paragraph 2 is the
    # blank line right after it. Together they are the leading 27
    # characters (26 for the import line incl. its paragraph mark, plus
    # 1 for the empty paragraph's mark).
    $lead = $tr.Characters(1, 27)
    $lead.Delete()
}
